# Rename the header columns to match the Jira export (add spaces) and
# drop the helper "Forecast"/"Done" ratio columns (I:L) that used to
# live next to the main table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sprint Name"
$ws.Range("B1").Value = "Start Date"
$ws.Range("C1").Value = "End Date"
$ws.Range("D1").Value = "Capacity Forecast"
$ws.Range("E1").Value = "Effort Forecast"
$ws.Range("F1").Value = "Capacity Done"
$ws.Range("G1").Value = "Effort Done"

# Remove the old helper formulas/values in columns I:L (rows 1-3 had content).
$ws.Range("I1:L3").ClearContents()

# Update Sprint 4 (row 5) effort/capacity-done values.
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 2

# Update Sprint 5 (row 6): end date, capacity forecast, effort forecast.
$ws.Range("C6").Value = 41763
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 5

# Update Sprint 6 (row 7): start/end dates and capacity forecast.
$ws.Range("B7").Value = 41764
$ws.Range("C7").Value = 41770
$ws.Range("D7").Value = 12

# Update Sprint 7 (row 8): start/end dates and capacity forecast.
$ws.Range("B8").Value = 41771
$ws.Range("C8").Value = 41777
$ws.Range("D8").Value = 12

# Update Sprint 8 (row 9): start/end dates and capacity forecast.
$ws.Range("B9").Value = 41778
$ws.Range("C9").Value = 41784
$ws.Range("D9").Value = 12

# Match the author's last on-screen selection.
$ws.Range("H1:L6").Select()
